$d = $word.ActiveDocument

# Position at the very end of the document body (after the last paragraph's
# text, before the final section break) so the new content is appended as
# new paragraphs following "Lesson-8-viewGroup-linearLayout".
$endRange = $d.Content
$endRange.Collapse(0)

# Build the OOXML fragment for the three new paragraphs. Using InsertXML lets
# us create paragraphs with multiple, distinctly-split <w:r> runs (matching
# the target markup) instead of having Word silently coalesce consecutive
# InsertAfter calls that share identical formatting into a single run.
$newParagraphsXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Layout_with/layout_height: thuộc tính</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Wrap</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>_</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>content</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>: bao nội dung bên trong</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Wrap</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>_</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>parent</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>: toàn bộ màn hình</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$endRange.InsertXML($newParagraphsXml)
